# random_items.xlsx survey tweak:
#  - bump the sampled-diamonds count from 175 to 180 (and tidy the call)
#  - swap the old readRDS(gzcon(url(...))) item-fetching calls for the
#    faster get_opencpu_rds() helper
#  - replace the long "personal training diary" note with a short dev note
#  - extend the rating-item block from q1..q12 up to q1..q22 (new rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- row 2 / row 3: the `item_list` / `items_seen` calculate rows ---------
$ws.Range("J2").Value = "library(ggplot2)`ndata(diamonds)`ndiamonds[sample(nrow(diamonds),size = 180) ,]"
$ws.Range("J3").Value = "as.character(jsonlite::toJSON(get_opencpu_rds(paste0(random_items`$item_list,""R/.val/rds""))[,'table']))"

# --- row 4: the intro note text --------------------------------------------
$ws.Range("F4").Value = "## Just testing how fast we can get David's items"
$ws.Rows.Item(4).RowHeight = 15

# --- rows 5..16 already hold q1..q12; just refresh the F-column formula ---
for ($i = 1; $i -le 12; $i++) {
    $row = $i + 4
    $ws.Cells.Item($row, 6).Value = "``r get_opencpu_rds(paste0(random_items`$item_list,""R/.val/rds""))[$i,'carat']``"
}

# --- rows 17..26: add q13..q22 rating items --------------------------------
for ($i = 13; $i -le 22; $i++) {
    $row = $i + 4
    $ws.Rows.Item($row).RowHeight = $(if ($row -eq 17) { 49 } else { 45 })
    $ws.Cells.Item($row, 3).Value = "rating_button 5"
    $ws.Cells.Item($row, 4).Value = "q$i"
    $ws.Cells.Item($row, 6).Value = "``r get_opencpu_rds(paste0(random_items`$item_list,""R/.val/rds""))[$i,'carat']``"
    $ws.Cells.Item($row, 7).Value = "low"
    $ws.Cells.Item($row, 8).Value = "high"
}
